# Correction liste distribution des cartes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments ---
$ws.Rows.Item(1).RowHeight = 27.7
$ws.Rows.Item(4).RowHeight = 20.25
$ws.Rows.Item(5).RowHeight = 20.05
$ws.Rows.Item(6).RowHeight = 32.05

# --- Row 4 : now holds "Plaque de tôle" set ---
$ws.Range("B4").Value = "Plaque de tôle"
$ws.Range("C4").Value = "Plaque de tôle"
$ws.Range("D4").Value = "Plaque de tôle"
$ws.Range("E4").Value = "Plaque de tôle"

# --- Row 5 : now holds "Sandwich / Bouteille d'eau" set ---
$ws.Range("B5").Value = "Sandwich"
$ws.Range("C5").Value = "Bouteille d’eau"
$ws.Range("D5").Value = "Sandwich"
$ws.Range("E5").Value = "Bouteille d’eau"

# --- Row 6 : now holds "Panier garni / Kit BBQ / Moulins" set ---
$ws.Range("B6").Value = "Panier garni"
$ws.Range("C6").Value = "Kit BBQ Cannibale"
$ws.Range("D6").Value = "Moulin à légumes"
$ws.Range("E6").Value = "Moulin à légumes inversé"

# --- Spelling / wording corrections ---
$ws.Range("E7").Value = "Hache"
$ws.Range("E8").Value = "Jeu de société Quoridor"

# --- Card re-distribution corrections (tours 3 à 6) ---
$ws.Range("E9").Value = "Cartouche"
$ws.Range("C10").Value = "Gourde"
$ws.Range("E10").Value = "Taser"
$ws.Range("D11").Value = "Magazine minceur"
$ws.Range("E11").Value = "Sandwich"
$ws.Range("C12").Value = "Somnifères"
$ws.Range("E12").Value = "Eau croupie"
